# Reproduce the "Add files via upload" commit:
#  - header cell A1 text changed from "PatientName" to "patient_name"
#  - column A given a custom width (~27.43 characters)
#  - active/selected cell moved to A3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header cell A1
$ws.Range("A1").Value = "patient_name"

# 2. Widen column A to fit the new header text (custom width)
$ws.Columns("A").ColumnWidth = 26.67

# 3. Leave A3 selected / active, matching the saved selection state
$ws.Range("A3").Select()
